$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-05-13 Monday" "2024-05-14 Tuesday"

Replace-Text "891÷2=" "120÷2="
Replace-Text "618÷5=" "828÷6="
Replace-Text "560÷6=" "702÷5="
Replace-Text "187÷5=" "984÷8="
Replace-Text "267÷6=" "951÷2="
Replace-Text "406÷3=" "391÷9="
Replace-Text "171÷7=" "415÷2="
Replace-Text "523÷6=" "437÷9="
Replace-Text "433÷8=" "707÷8="
Replace-Text "929÷3=" "340÷3="
Replace-Text "527÷9=" "808÷3="
Replace-Text "751÷3=" "511÷8="
Replace-Text "874÷8=" "248÷9="
Replace-Text "658÷7=" "734÷3="
Replace-Text "646÷8=" "502÷5="
Replace-Text "202÷4=" "857÷3="
Replace-Text "138÷3=" "495÷8="
Replace-Text "577÷2=" "179÷4="
Replace-Text "129÷2=" "870÷7="
Replace-Text "563÷5=" "749÷3="
Replace-Text "930÷8=" "921÷3="
Replace-Text "661÷2=" "855÷3="
Replace-Text "314÷8=" "772÷2="
Replace-Text "385÷6=" "376÷6="
Replace-Text "649÷7=" "506÷6="
